# Groups1.xlsx — refresh the carpool roster with a new (smaller) set of
# children and re-pack the trailing "school"/"cost"/"time" summary rows
# directly beneath them. (Point colors are matched to the other graphs
# outside this sheet, per the commit message — only the data changes here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cells whose new text would otherwise be auto-coerced to a Number
#     (losing a literal trailing ".0") need to be pre-formatted as Text.
#     (Union ranges only apply NumberFormat to their first area, so set
#     each cell individually.) ---
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H8").NumberFormat = "@"
$ws.Range("H9").NumberFormat = "@"
$ws.Range("H10").NumberFormat = "@"
$ws.Range("B13").NumberFormat = "@"

# nChildren
$ws.Range("B4").Value = "5"

# Row 6
$ws.Range("B6").Value = "8"
$ws.Range("C6").Value = "Marni  "
$ws.Range("D6").Value = "Shanika  "
$ws.Range("E6").Value = "3.52,-5.3"
$ws.Range("F6").Value = "Lady(mother): 0560804012"
$ws.Range("H6").Value = "31.0"

# Row 7
$ws.Range("B7").Value = "7"
$ws.Range("C7").Value = "Wyatt  "
$ws.Range("D7").Value = "Willette  "
$ws.Range("E7").Value = "5.19,-5.63"
$ws.Range("F7").Value = "Antionette(father): 0557331799"
$ws.Range("G7").Value = "7:03:00"
$ws.Range("H7").Value = "28.0"

# Row 8
$ws.Range("B8").Value = "1"
$ws.Range("C8").Value = "Corene  "
$ws.Range("D8").Value = "Myra  "
$ws.Range("E8").Value = "7.75,-1.13"
$ws.Range("F8").Value = "Georgie(mother): 0544823581"
$ws.Range("G8").Value = "7:09:00"
$ws.Range("H8").Value = "22.0"

# Row 9
$ws.Range("B9").Value = "6"
$ws.Range("C9").Value = "Ema  "
$ws.Range("D9").Value = "Ardell  "
$ws.Range("E9").Value = "6.04,5.4"
$ws.Range("F9").Value = "Carley(grandmother): 0533587167"
$ws.Range("G9").Value = "7:18:00"
$ws.Range("H9").Value = "13.0"

# Row 10
$ws.Range("B10").Value = "13"
$ws.Range("C10").Value = "Fay  "
$ws.Range("D10").Value = "Emilee  "
$ws.Range("E10").Value = "4.45,1.94"
$ws.Range("F10").Value = "Sheri(mother): 0516797453"
$ws.Range("G10").Value = "7:24:00"
$ws.Range("H10").Value = "7.0"

# Row 11 becomes the "school" trailer row (used to be a child row); it has
# no H value, so clear what used to be there.
$ws.Range("A11").Value = "school"
$ws.Range("B11").Value = "3"
$ws.Range("C11").Value = "Ironiah"
$ws.Range("D11").Value = "mySchool"
$ws.Range("E11").Value = "0,0"
$ws.Range("F11").Value = "Shir(secretary): 0523345098"
$ws.Range("G11").Value = "7:31:00"
$ws.Range("H11").ClearContents()

# Row 12 becomes the "cost" trailer row; drop the old child columns C:H.
$ws.Range("A12").Value = "cost"
$ws.Range("B12").Value = "25"
$ws.Range("C12:H12").ClearContents()

# Row 13 becomes the "time" trailer row; drop the old child columns C:H.
$ws.Range("A13").Value = "time"
$ws.Range("B13").Value = "31.0"
$ws.Range("C13:H13").ClearContents()

# Old rows 14-16 (previous school/cost/time trailer) are no longer needed.
$ws.Range("A14:H16").ClearContents()
